# Apply cryptos list update (Mon Dec 11 17:26:15 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.685.44'
$ws.Range("E2").Value = '  -4.85%  '
$ws.Range("D3").Value = '2.206.88'
$ws.Range("E3").Value = '  -6.01%  '
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("D5").Value = '245.17'
$ws.Range("E5").Value = '  +2.42%  '
$ws.Range("D6").Value = '0.626'
$ws.Range("E6").Value = '  -5.85%  '
$ws.Range("D7").Value = '70.34'
$ws.Range("E7").Value = '  -3.84%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").Value = '0.546'
$ws.Range("E9").Value = '  -8.42%  '
$ws.Range("D10").Value = '36.91'
$ws.Range("E10").Value = '  +9.96%  '
$ws.Range("D11").Value = '0.0945'
$ws.Range("E11").Value = '  -6.29%  '
$ws.Range("D12").Value = '57.99'
$ws.Range("E12").Value = '  -5.87%  '
$ws.Range("E13").Value = '  -3.02%  '
$ws.Range("D14").Value = '6.67'
$ws.Range("E14").Value = '  -7.42%  '
$ws.Range("D15").Value = '2.538.42'
$ws.Range("E15").Value = '  -5.93%  '
$ws.Range("D16").Value = '14.78'
$ws.Range("E16").Value = '  -8.29%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.840'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -6.93%  '
$ws.Range("D18").Value = '2.204.03'
$ws.Range("E18").Value = '  -6.06%  '
$ws.Range("D19").Value = '41.614.89'
$ws.Range("E19").Value = '  -4.94%  '
$ws.Range("D20").Value = '0.0₃0952'
$ws.Range("E20").Value = '  -6.99%  '
$ws.Range("D21").Value = '73.41'
$ws.Range("E21").Value = '  -5.68%  '
$ws.Range("D22").Value = '6.06'
$ws.Range("E22").Value = '  -7.66%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '233.90'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -7.36%  '
$ws.Range("D24").Value = '2.06'
$ws.Range("E24").Value = '  +11.77%  '
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  -0.08%  '
$ws.Range("E26").Value = '  -5.56%  '
$ws.Range("D27").Value = '2.44'
$ws.Range("E27").Value = '  -1.70%  '
$ws.Range("D28").Value = '2.22'
$ws.Range("E28").Value = '  -2.08%  '
$ws.Range("D29").Value = '9.84'
$ws.Range("E29").Value = '  -5.36%  '
$ws.Range("D30").Value = '169.36'
$ws.Range("E30").Value = '  -3.59%  '
$ws.Range("D31").Value = '20.38'
$ws.Range("E32").Value = '  -6.40%  '
$ws.Range("D33").Value = '0.124'
$ws.Range("E33").Value = '  -7.25%  '
$ws.Range("D34").Value = '0.0711'
$ws.Range("E34").Value = '  -3.67%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.10'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.44%  '
$ws.Range("D36").Value = '4.59'
$ws.Range("E36").Value = '  -8.87%  '
$ws.Range("D37").Value = '3.87'
$ws.Range("E37").Value = '  +3.38%  '
$ws.Range("D38").Value = '23.37'
$ws.Range("E38").Value = '  +20.17%  '
$ws.Range("E39").Value = '  -5.11%  '
$ws.Range("D40").Value = '0.0272'
$ws.Range("E40").Value = '  +0.35%  '
$ws.Range("D41").Value = '5.84'
$ws.Range("E41").Value = '  -8.85%  '
$ws.Range("D42").Value = '65.24'
$ws.Range("E42").Value = '  +0.70%  '
$ws.Range("D43").Value = '8.94'
$ws.Range("E43").Value = '  -1.63%  '
$ws.Range("D44").Value = '4.87'
$ws.Range("E44").Value = '  -9.93%  '
$ws.Range("E45").Value = '  -3.39%  '
$ws.Range("D46").Value = '0.0996'
$ws.Range("E46").Value = '  -6.27%  '
$ws.Range("E47").Value = '  -0.05%  '
$ws.Range("D48").Value = '4.54'
$ws.Range("E48").Value = '  +5.67%  '
$ws.Range("D49").Value = '10.39'
$ws.Range("E49").Value = '  +9.33%  '
$ws.Range("E50").Value = '  -3.76%  '
$ws.Range("B51").Value = 'BitTorrent-New'
$ws.Range("C51").Value = 'https://coinranking.com/coin/w4MqH_Xe8+bittorrent-new-btt'
$ws.Range("D51").Value = '0.0₃0148'
$ws.Range("E51").Value = '  +10.60%  '
